$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = 81
$ws.Range("B82").Value = 1
$ws.Range("C82").Value = "2024-06-16 15:13:21"
$ws.Range("D82").Value = 200
$ws.Range("E82").Value = 6

$ws.Range("A83").Value = 82
$ws.Range("B83").Value = 2
$ws.Range("C83").Value = "2024-06-16 15:13:22"
$ws.Range("D83").Value = 200
$ws.Range("E83").Value = 1
